$d = $word.ActiveDocument

# Reusable OOXML wrapper for Range.InsertXML() snippets.
function New-WordPackageXml($innerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $innerBodyXml + '</w:body>' + `
        '</w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Drop the stray "_GoBack" bookmark that used to sit right after the
#    "Space Borders" title run.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) "Shift to inner borders layer" -> "Move" + " to inner borders layer"
#    (split into two separate runs, first run's text literally changes
#    from "Shift" to "Move").
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Shift to inner borders layer*") {
        $r = $p.Range
        # Wipe the paragraph's text (keep the paragraph mark) then rebuild
        # it from literal OOXML so the two pieces land in distinct <w:r>s.
        $body = $d.Range($r.Start, $r.End - 1)
        $body.Text = ""
        $insertAt = $d.Range($r.Start, $r.Start)
        $runsXml = '<w:p>' + `
            '<w:r><w:t>Move</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve"> to inner borders layer</w:t></w:r>' + `
            '</w:p>'
        $insertAt.InsertXML((New-WordPackageXml $runsXml))
        break
    }
}

# ---------------------------------------------------------------------
# 3) Append extra sentences (and the re-inserted "_GoBack" bookmark) to
#    the very last paragraph, right after its existing " " run.
# ---------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastRange = $lastPara.Range
$wholePara = $d.Range($lastRange.Start, $lastRange.End)
$newParaXml = '<w:p>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>For expanding borders, just do the same but use the regular Brush tool</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> instead of the Pencil,</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t xml:space="preserve"> with hard brush selected, size 9px. </w:t></w:r>' + `
    '</w:p>'
$wholePara.InsertXML((New-WordPackageXml $newParaXml))

Write-Output "Edit applied."
